# Update the report folder path used in the workbook's stored absolute path
# (June -> July) is metadata Excel regenerates from the save location; the
# meaningful, scriptable change is the new data row for July.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Append the new "Jul, 2021" data row right after the last existing row (31)
$ws.Range("A32").Value = "Jul, 2021"
$ws.Range("B32").Value = 501
$ws.Range("C32").Value = 149

# Match the author's final selection left after entering the new row
$ws.Range("D28").Select()
